# Auto-generated edit script: updates profit-calculation columns (H:N)
# across multiple job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to
# reflect refreshed market-price data from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value2 = 648.4
$ws.Range("I2").Value2 = 572.125
$ws.Range("J2").Value2 = 735.5714
$ws.Range("K2").Value2 = 572.125
$ws.Range("L2").Value2 = 735.5714
$ws.Range("M2").Value2 = -459.125
$ws.Range("N2").Value2 = -961.5714

$ws.Range("H103").Value2 = 602.551
$ws.Range("I103").Value2 = 512.7907
$ws.Range("J103").Value2 = 1245.8334
$ws.Range("K103").Value2 = 1538.3721
$ws.Range("L103").Value2 = 3737.5002
$ws.Range("M103").Value2 = -952.3721
$ws.Range("N103").Value2 = -4909.5002

$ws.Range("H123").Value2 = 71999
$ws.Range("J123").Value2 = 71999
$ws.Range("L123").Value2 = 71999
$ws.Range("N123").Value2 = -81799

$ws.Range("H131").Value2 = 15173.182
$ws.Range("I131").Value2 = 1364.5
$ws.Range("J131").Value2 = 51996.332
$ws.Range("K131").Value2 = 4093.5
$ws.Range("L131").Value2 = 155988.996
$ws.Range("M131").Value2 = 946.5
$ws.Range("N131").Value2 = -166068.996

$ws.Range("H135").Value2 = 1023.21875
$ws.Range("I135").Value2 = 632.9583
$ws.Range("J135").Value2 = 2194
$ws.Range("K135").Value2 = 5696.6247
$ws.Range("L135").Value2 = 19746
$ws.Range("M135").Value2 = -3161.6247
$ws.Range("N135").Value2 = -24816

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value2 = 167308.17
$ws.Range("I16").Value2 = 1000000
$ws.Range("K16").Value2 = 1000000
$ws.Range("M16").Value2 = -999713

$ws.Range("H32").Value2 = 302007.72
$ws.Range("I32").Value2 = 340608.38
$ws.Range("K32").Value2 = 340608.38
$ws.Range("M32").Value2 = -340321.38

$ws.Range("H45").Value2 = 3506.2354
$ws.Range("I45").Value2 = 2900.5
$ws.Range("K45").Value2 = 2900.5
$ws.Range("M45").Value2 = -2523.5

$ws.Range("H61").Value2 = 5497.364
$ws.Range("I61").Value2 = 1972.6428
$ws.Range("K61").Value2 = 1972.6428
$ws.Range("M61").Value2 = -1760.6428

$ws.Range("H132").Value2 = 2506992
$ws.Range("I132").Value2 = 3576287.2
$ws.Range("K132").Value2 = 10728861.6
$ws.Range("M132").Value2 = -10726331.6

$ws.Range("H134").Value2 = 64994
$ws.Range("J134").Value2 = 64994
$ws.Range("L134").Value2 = 64994
$ws.Range("N134").Value2 = -75134

$ws.Range("H136").Value2 = 5497.364
$ws.Range("I136").Value2 = 1972.6428
$ws.Range("K136").Value2 = 5917.928400000001
$ws.Range("M136").Value2 = -3367.928400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value2 = 3399.9092
$ws.Range("I86").Value2 = 3880
$ws.Range("K86").Value2 = 3880
$ws.Range("M86").Value2 = -2757

$ws.Range("H89").Value2 = 3399.9092
$ws.Range("I89").Value2 = 3880
$ws.Range("K89").Value2 = 19400
$ws.Range("M89").Value2 = -13784

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 2162.0186
$ws.Range("I31").Value2 = 2094.1936
$ws.Range("K31").Value2 = 2094.1936
$ws.Range("M31").Value2 = -1799.1936

$ws.Range("H34").Value2 = 2162.0186
$ws.Range("I34").Value2 = 2094.1936
$ws.Range("K34").Value2 = 2094.1936
$ws.Range("M34").Value2 = -1892.1936

$ws.Range("H58").Value2 = 6439.1934
$ws.Range("I58").Value2 = 3304.647
$ws.Range("K58").Value2 = 3304.647
$ws.Range("M58").Value2 = -3101.647

$ws.Range("H94").Value2 = 5903.8184
$ws.Range("I94").Value2 = 11325.4
$ws.Range("K94").Value2 = 11325.4
$ws.Range("M94").Value2 = -10874.4

$ws.Range("H107").Value2 = 749.8182
$ws.Range("I107").Value2 = 694.3333
$ws.Range("K107").Value2 = 694.3333
$ws.Range("M107").Value2 = 1225.6667

$ws.Range("H122").Value2 = 5869.825
$ws.Range("I122").Value2 = 1519.2858
$ws.Range("J122").Value2 = 16021.083
$ws.Range("K122").Value2 = 4557.857400000001
$ws.Range("L122").Value2 = 48063.249
$ws.Range("M122").Value2 = -2107.857400000001
$ws.Range("N122").Value2 = -52963.249

$ws.Range("H136").Value2 = 6439.1934
$ws.Range("I136").Value2 = 3304.647
$ws.Range("K136").Value2 = 9913.940999999999
$ws.Range("M136").Value2 = -7363.940999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value2 = 58.3
$ws.Range("I2").Value2 = 56.75
$ws.Range("J2").Value2 = 59.333332
$ws.Range("K2").Value2 = 340.5
$ws.Range("L2").Value2 = 355.999992
$ws.Range("M2").Value2 = -227.5
$ws.Range("N2").Value2 = -581.999992

$ws.Range("H103").Value2 = 377.57144
$ws.Range("I103").Value2 = 164.66667
$ws.Range("J103").Value2 = 537.25
$ws.Range("K103").Value2 = 494.00001
$ws.Range("L103").Value2 = 1611.75
$ws.Range("M103").Value2 = 384.99999
$ws.Range("N103").Value2 = -3369.75

$ws.Range("H113").Value2 = 1154.5
$ws.Range("I113").Value2 = 623.625
$ws.Range("K113").Value2 = 1870.875
$ws.Range("M113").Value2 = 299.125

$ws.Range("H128").Value2 = 0
$ws.Range("I128").Value2 = 0
$ws.Range("K128").Value2 = 0
$ws.Range("M128").ClearContents()

$ws.Range("H131").Value2 = 2675.1892
$ws.Range("I131").Value2 = 805
$ws.Range("J131").Value2 = 2782.0571
$ws.Range("K131").Value2 = 2415
$ws.Range("L131").Value2 = 8346.1713
$ws.Range("M131").Value2 = 2625
$ws.Range("N131").Value2 = -18426.1713

$ws.Range("H132").Value2 = 828.5833
$ws.Range("J132").Value2 = 1150
$ws.Range("L132").Value2 = 10350
$ws.Range("N132").Value2 = -15410

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value2 = 1433333.4
$ws.Range("I11").Value2 = 1900000
$ws.Range("K11").Value2 = 1900000
$ws.Range("M11").Value2 = -1899861

$ws.Range("H14").Value2 = 66735930
$ws.Range("I14").Value2 = 71502280
$ws.Range("K14").Value2 = 71502280
$ws.Range("M14").Value2 = -71502112

$ws.Range("H18").Value2 = 13004.2
$ws.Range("I18").Value2 = 13004.2
$ws.Range("K18").Value2 = 13004.2
$ws.Range("M18").Value2 = -12711.2

$ws.Range("H110").Value2 = 52000
$ws.Range("J110").Value2 = 52000
$ws.Range("L110").Value2 = 52000
$ws.Range("N110").Value2 = -60180

$ws.Range("H122").Value2 = 170657.67
$ws.Range("I122").Value2 = 253361.75
$ws.Range("K122").Value2 = 760085.25
$ws.Range("M122").Value2 = -757635.25

$ws.Range("H126").Value2 = 2467.9167
$ws.Range("I126").Value2 = 2161.353
$ws.Range("K126").Value2 = 6484.059
$ws.Range("M126").Value2 = -4014.059

$ws.Range("H132").Value2 = 13272.391
$ws.Range("I132").Value2 = 17242.178
$ws.Range("K132").Value2 = 51726.534
$ws.Range("M132").Value2 = -49196.534

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 1873.3334
$ws.Range("I22").Value2 = 849.5
$ws.Range("J22").Value2 = 1933.5588
$ws.Range("K22").Value2 = 849.5
$ws.Range("L22").Value2 = 1933.5588
$ws.Range("M22").Value2 = -554.5
$ws.Range("N22").Value2 = -2523.5588

$ws.Range("H23").Value2 = 5000
$ws.Range("I23").Value2 = 0
$ws.Range("J23").Value2 = 5000
$ws.Range("K23").Value2 = 0
$ws.Range("L23").Value2 = 5000
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value2 = -5460

$ws.Range("H25").Value2 = 0
$ws.Range("J25").Value2 = 0
$ws.Range("L25").Value2 = 0
$ws.Range("N25").ClearContents()

$ws.Range("H27").Value2 = 1873.3334
$ws.Range("I27").Value2 = 849.5
$ws.Range("J27").Value2 = 1933.5588
$ws.Range("K27").Value2 = 849.5
$ws.Range("L27").Value2 = 1933.5588
$ws.Range("M27").Value2 = -742.5
$ws.Range("N27").Value2 = -2147.5588

$ws.Range("H93").Value2 = 3907.2856
$ws.Range("I93").Value2 = 1482.8
$ws.Range("K93").Value2 = 1482.8
$ws.Range("M93").Value2 = -234.8

$ws.Range("H136").Value2 = 9940.777
$ws.Range("I136").Value2 = 5552
$ws.Range("K136").Value2 = 16656
$ws.Range("M136").Value2 = -14106

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value2 = 3315.6843
$ws.Range("I136").Value2 = 3066.6667
$ws.Range("J136").Value2 = 4249.5
$ws.Range("K136").Value2 = 9200.000100000001
$ws.Range("L136").Value2 = 12748.5
$ws.Range("M136").Value2 = -6650.000100000001
$ws.Range("N136").Value2 = -17848.5
